$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "2565"

# Delete column A (the numeric index column), shifting B->A, C->B, D->C
$ws.Columns.Item(1).Delete()

# Remove all cell styling (bold header font, borders, alignment) so cells use the default style
$ws.UsedRange.ClearFormats()

# Clear the now-empty "Face Status" data cells (C2:C21) so no cell objects remain there
$ws.Range("C2:C21").Clear()
